$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Add materials for session 04 (row 5): task link (Aufgaben) then slides link (Folien)
$ws.Range("F5").Value = "https://stats.ifp.uni-mainz.de/ba-ccs-track/befragung-appkit.html"
$ws.Range("E5").Value = "slides/slides.html#/sitzung-04-situative-befragungsdesigns"

# Update the active selection to E5
$ws.Range("E5").Select()
